$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{Col="D"; Row=2; Value="66.790.97"},
    @{Col="E"; Row=2; Value="  +1.72%  "},
    @{Col="D"; Row=3; Value="3.797.54"},
    @{Col="E"; Row=3; Value="  +2.21%  "},
    @{Col="E"; Row=4; Value="  +0.43%  "},
    @{Col="D"; Row=5; Value="421.10"},
    @{Col="E"; Row=5; Value="  +2.49%  "},
    @{Col="D"; Row=6; Value="129.13"},
    @{Col="E"; Row=6; Value="  -4.77%  "},
    @{Col="D"; Row=7; Value="3.793.90"},
    @{Col="E"; Row=7; Value="  +2.67%  "},
    @{Col="D"; Row=8; Value="0.603"},
    @{Col="E"; Row=8; Value="  -3.70%  "},
    @{Col="E"; Row=9; Value="  -0.18%  "},
    @{Col="D"; Row=10; Value="0.721"},
    @{Col="E"; Row=10; Value="  -2.07%  "},
    @{Col="E"; Row=11; Value="  -2.57%  "},
    @{Col="D"; Row=12; Value="0.0000352"},
    @{Col="E"; Row=12; Value="  +6.29%  "},
    @{Col="D"; Row=13; Value="40.05"},
    @{Col="E"; Row=13; Value="  -6.04%  "},
    @{Col="B"; Row=14; Value="Polkadot"},
    @{Col="C"; Row=14; Value="https://coinranking.com/coin/25W7FG7om+polkadot-dot"},
    @{Col="D"; Row=14; Value="10.09"},
    @{Col="E"; Row=14; Value="  +0.07%  "},
    @{Col="B"; Row=15; Value="WrappedliquidstakedEther2.0"},
    @{Col="C"; Row=15; Value="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"},
    @{Col="D"; Row=15; Value="4.386.42"},
    @{Col="E"; Row=15; Value="  +1.87%  "},
    @{Col="D"; Row=16; Value="15.70"},
    @{Col="E"; Row=16; Value="  +20.38%  "},
    @{Col="E"; Row=17; Value="  -1.00%  "},
    @{Col="D"; Row=18; Value="3.783.24"},
    @{Col="E"; Row=18; Value="  +2.03%  "},
    @{Col="D"; Row=19; Value="19.47"},
    @{Col="E"; Row=19; Value="  -3.28%  "},
    @{Col="D"; Row=20; Value="66.844.42"},
    @{Col="E"; Row=20; Value="  +2.06%  "},
    @{Col="E"; Row=21; Value="  -2.30%  "},
    @{Col="D"; Row=22; Value="406.56"},
    @{Col="E"; Row=22; Value="  -4.33%  "},
    @{Col="D"; Row=23; Value="14.25"},
    @{Col="E"; Row=23; Value="  -3.00%  "},
    @{Col="D"; Row=24; Value="83.66"},
    @{Col="E"; Row=24; Value="  -3.85%  "},
    @{Col="D"; Row=25; Value="3.00"},
    @{Col="E"; Row=25; Value="  -0.72%  "},
    @{Col="D"; Row=26; Value="37.03"},
    @{Col="E"; Row=26; Value="  +1.40%  "},
    @{Col="D"; Row=27; Value="5.58"},
    @{Col="E"; Row=27; Value="  +8.46%  "},
    @{Col="D"; Row=28; Value="3.19"},
    @{Col="E"; Row=28; Value="  +0.02%  "},
    @{Col="D"; Row=29; Value="9.38"},
    @{Col="E"; Row=29; Value="  -2.02%  "},
    @{Col="D"; Row=30; Value="713.03"},
    @{Col="E"; Row=30; Value="  +5.51%  "},
    @{Col="D"; Row=31; Value="8.37"},
    @{Col="E"; Row=31; Value="  +19.38%  "},
    @{Col="B"; Row=32; Value="Cosmos"},
    @{Col="C"; Row=32; Value="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"},
    @{Col="D"; Row=32; Value="12.39"},
    @{Col="E"; Row=32; Value="  -1.18%  "},
    @{Col="B"; Row=33; Value="Toncoin"},
    @{Col="C"; Row=33; Value="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"},
    @{Col="D"; Row=33; Value="2.77"},
    @{Col="E"; Row=33; Value="  +1.68%  "},
    @{Col="E"; Row=34; Value="  +0.54%  "},
    @{Col="D"; Row=35; Value="0.999"},
    @{Col="E"; Row=35; Value="  +0.02%  "},
    @{Col="E"; Row=36; Value="  -5.37%  "},
    @{Col="D"; Row=37; Value="38.41"},
    @{Col="E"; Row=37; Value="  -6.64%  "},
    @{Col="D"; Row=38; Value="54.95"},
    @{Col="E"; Row=38; Value="  -1.62%  "},
    @{Col="D"; Row=39; Value="0.0₃0758"},
    @{Col="E"; Row=39; Value="  +12.06%  "},
    @{Col="B"; Row=40; Value="NEARProtocol"},
    @{Col="C"; Row=40; Value="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"},
    @{Col="D"; Row=40; Value="4.98"},
    @{Col="E"; Row=40; Value="  +16.83%  "},
    @{Col="B"; Row=41; Value="VeChain"},
    @{Col="C"; Row=41; Value="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"},
    @{Col="D"; Row=41; Value="0.0450"},
    @{Col="E"; Row=41; Value="  -4.10%  "},
    @{Col="D"; Row=42; Value="2.94"},
    @{Col="E"; Row=42; Value="  -0.04%  "},
    @{Col="E"; Row=43; Value="  +1.08%  "},
    @{Col="E"; Row=44; Value="  -5.03%  "},
    @{Col="D"; Row=45; Value="3.33"},
    @{Col="E"; Row=45; Value="  -0.43%  "},
    @{Col="D"; Row=46; Value="143.94"},
    @{Col="E"; Row=46; Value="  -0.67%  "},
    @{Col="D"; Row=47; Value="3.11"},
    @{Col="E"; Row=47; Value="  +0.34%  "},
    @{Col="D"; Row=48; Value="2.03"},
    @{Col="E"; Row=48; Value="  -1.99%  "},
    @{Col="D"; Row=49; Value="25.60"},
    @{Col="E"; Row=49; Value="  +1.28%  "},
    @{Col="B"; Row=50; Value="Stacks"},
    @{Col="C"; Row=50; Value="https://coinranking.com/coin/mMPrMcB7+stacks-stx"},
    @{Col="D"; Row=50; Value="2.78"},
    @{Col="E"; Row=50; Value="  -1.18%  "},
    @{Col="B"; Row=51; Value="WEMIXToken"},
    @{Col="C"; Row=51; Value="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"},
    @{Col="D"; Row=51; Value="2.54"},
    @{Col="E"; Row=51; Value="  +0.86%  "}
)

foreach ($e in $edits) {
    $cell = $ws.Range($e.Col + $e.Row)
    $cell.Value = "'" + $e.Value
    $cell.Style = $ws.Range("B" + $e.Row).Style
}